$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, shifting existing rows 8-28 down to 9-29
$ws.Rows(8).Insert()

# Populate the newly inserted row 8 (copy of what is now row 9, with Fecha and Volumen updated)
$ws.Cells.Item(8, 1).Value = 10
$ws.Cells.Item(8, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(8, 3).Value = "La Araucanía"
$ws.Cells.Item(8, 4).Value = 44469
$ws.Cells.Item(8, 4).Style = $ws.Cells.Item(9, 4).Style
$ws.Cells.Item(8, 4).NumberFormat = $ws.Cells.Item(9, 4).NumberFormat
$ws.Cells.Item(8, 5).Value = 9
$ws.Cells.Item(8, 6).Value = 100112026
$ws.Cells.Item(8, 7).Value = "Haba"
$ws.Cells.Item(8, 8).Value = "Sin especificar"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 90
$ws.Cells.Item(8, 11).Value = 14000
$ws.Cells.Item(8, 12).Value = 14000
$ws.Cells.Item(8, 13).Value = 14000
$ws.Cells.Item(8, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(8, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(8, 16).Value = 560
$ws.Cells.Item(8, 17).Value = 25
$ws.Cells.Item(8, 18).Value = "Hortaliza"
